# Scheduled-runner refresh of the cached market-price / profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) across the
# eight crafter-job sheets. These are plain cached values (no formulas in
# the sheet), so each refreshed figure is written straight to its cell;
# a couple of rows had an M/N value newly appear (now computable) or
# disappear (no longer applicable), handled with ClearContents().
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 566.6667
$ws.Range("I9").Value = 66.666664
$ws.Range("J9").Value = 1066.6666
$ws.Range("K9").Value = 66.666664
$ws.Range("L9").Value = 1066.6666
$ws.Range("M9").Value = 102.333336
$ws.Range("N9").Value = -1404.6666
$ws.Range("H29").Value = 4000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H33").Value = 318.72415
$ws.Range("J33").Value = 326.69232
$ws.Range("L33").Value = 326.69232
$ws.Range("N33").Value = -784.69232
$ws.Range("H88").Value = 7969.4443
$ws.Range("I88").Value = 7980
$ws.Range("J88").Value = 7962.727
$ws.Range("K88").Value = 7980
$ws.Range("L88").Value = 7962.727
$ws.Range("M88").Value = -7574
$ws.Range("N88").Value = -8774.726999999999
$ws.Range("H91").Value = 7969.4443
$ws.Range("I91").Value = 7980
$ws.Range("J91").Value = 7962.727
$ws.Range("K91").Value = 7980
$ws.Range("L91").Value = 7962.727
$ws.Range("M91").Value = -6576
$ws.Range("N91").Value = -10770.727
$ws.Range("H98").Value = 4176.1333
$ws.Range("I98").Value = 5012.909
$ws.Range("K98").Value = 5012.909
$ws.Range("M98").Value = -3514.909
$ws.Range("H113").Value = 3490.4
$ws.Range("I113").Value = 4027
$ws.Range("J113").Value = 2417.2
$ws.Range("K113").Value = 4027
$ws.Range("L113").Value = 2417.2
$ws.Range("M113").Value = -773
$ws.Range("N113").Value = -8925.200000000001
$ws.Range("H122").Value = 4176.1333
$ws.Range("I122").Value = 5012.909
$ws.Range("K122").Value = 15038.727
$ws.Range("M122").Value = -12588.727
$ws.Range("H138").Value = 1703.9727
$ws.Range("I138").Value = 1008
$ws.Range("J138").Value = 3548.3
$ws.Range("K138").Value = 3024
$ws.Range("L138").Value = 10644.9
$ws.Range("M138").Value = 2116
$ws.Range("N138").Value = -20924.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1382.5555
$ws.Range("I74").Value = 1430.375
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 1430.375
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -556.375
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 1382.5555
$ws.Range("I77").Value = 1430.375
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 7151.875
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -2783.875
$ws.Range("N77").Value = -13736
$ws.Range("H122").Value = 2382.8125
$ws.Range("I122").Value = 2268.5
$ws.Range("K122").Value = 6805.5
$ws.Range("M122").Value = -4355.5
$ws.Range("H132").Value = 1600.0834
$ws.Range("I132").Value = 1414.6
$ws.Range("K132").Value = 4243.799999999999
$ws.Range("M132").Value = -1713.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3962.375
$ws.Range("I31").Value = 2762.6667
$ws.Range("J31").Value = 6252.727
$ws.Range("K31").Value = 2762.6667
$ws.Range("L31").Value = 6252.727
$ws.Range("M31").Value = -2467.6667
$ws.Range("N31").Value = -6842.727
$ws.Range("H34").Value = 3962.375
$ws.Range("I34").Value = 2762.6667
$ws.Range("J34").Value = 6252.727
$ws.Range("K34").Value = 2762.6667
$ws.Range("L34").Value = 6252.727
$ws.Range("M34").Value = -2560.6667
$ws.Range("N34").Value = -6656.727
$ws.Range("H37").Value = 10250
$ws.Range("I37").Value = 2000
$ws.Range("J37").Value = 13000
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 13000
$ws.Range("M37").Value = -1893
$ws.Range("N37").Value = -13214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 4000
$ws.Range("J57").Value = 4000
$ws.Range("L57").Value = 12000
$ws.Range("N57").Value = -13118
$ws.Range("H116").Value = 2064.2856
$ws.Range("J116").Value = 3750
$ws.Range("L116").Value = 11250
$ws.Range("N116").Value = -18134
$ws.Range("H118").Value = 1550
$ws.Range("I118").Value = 500
$ws.Range("J118").Value = 2600
$ws.Range("K118").Value = 1500
$ws.Range("L118").Value = 7800
$ws.Range("M118").Value = -257
$ws.Range("N118").Value = -10286
$ws.Range("H131").Value = 1697
$ws.Range("J131").Value = 1209.091
$ws.Range("L131").Value = 3627.273
$ws.Range("N131").Value = -13707.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 41261.25
$ws.Range("J48").Value = 41261.25
$ws.Range("L48").Value = 41261.25
$ws.Range("N48").Value = -42231.25
$ws.Range("H126").Value = 3255.04
$ws.Range("I126").Value = 2593.3333
$ws.Range("J126").Value = 3345.2727
$ws.Range("K126").Value = 7779.999899999999
$ws.Range("L126").Value = 10035.8181
$ws.Range("M126").Value = -5309.999899999999
$ws.Range("N126").Value = -14975.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1605.0555
$ws.Range("I7").Value = 1289.2222
$ws.Range("J7").Value = 1920.8889
$ws.Range("K7").Value = 1289.2222
$ws.Range("L7").Value = 1920.8889
$ws.Range("M7").Value = -1177.2222
$ws.Range("N7").Value = -2144.8889
$ws.Range("H126").Value = 1605.0555
$ws.Range("I126").Value = 1289.2222
$ws.Range("J126").Value = 1920.8889
$ws.Range("K126").Value = 3867.6666
$ws.Range("L126").Value = 5762.6667
$ws.Range("M126").Value = -1397.6666
$ws.Range("N126").Value = -10702.6667
$ws.Range("H128").Value = 29966.334
$ws.Range("J128").Value = 29966.334
$ws.Range("L128").Value = 29966.334
$ws.Range("N128").Value = -39926.334
$ws.Range("H132").Value = 5250.026
$ws.Range("I132").Value = 2043.8823
$ws.Range("J132").Value = 7727.5
$ws.Range("K132").Value = 6131.6469
$ws.Range("L132").Value = 23182.5
$ws.Range("M132").Value = -3601.6469
$ws.Range("N132").Value = -28242.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15298.929
$ws.Range("J54").Value = 15298.929
$ws.Range("L54").Value = 15298.929
$ws.Range("N54").Value = -16338.929
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H126").Value = 1946.8649
$ws.Range("I126").Value = 1701.44
$ws.Range("J126").Value = 2458.1667
$ws.Range("K126").Value = 5104.32
$ws.Range("L126").Value = 7374.500100000001
$ws.Range("M126").Value = -2634.32
$ws.Range("N126").Value = -12314.5001
$ws.Range("H132").Value = 18517.938
$ws.Range("I132").Value = 2821.182
$ws.Range("J132").Value = 53050.8
$ws.Range("K132").Value = 8463.545999999998
$ws.Range("L132").Value = 159152.4
$ws.Range("M132").Value = -5933.545999999998
$ws.Range("N132").Value = -164212.4
